$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ambush" (D5) has been folded into a renamed "Shocking Arrival" (D4),
# which is now called "Sudden Ambush". "Multiward" moves from E3 down to E4.
$ws.Range("D4").Value = "Sudden Ambush"
$ws.Range("E3").ClearContents()
$ws.Range("E4").Value = "Multiward"
$ws.Range("D5").ClearContents()

# Update the active selection to match the saved view state.
$ws.Range("E4").Select()
